# Update "想去人数" (want-to-go count) figures in the F column on both the
# "展览" and "全部类型" worksheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of sheet name -> hashtable of row number -> new F-column value
$updates = @{
    "展览" = @{
        4  = 1243
        10 = 3456
        12 = 85
        16 = 592
        17 = 86
        18 = 730
        19 = 207
        24 = 2602
        25 = 5107
        29 = 1533
        31 = 2239
        35 = 113
        37 = 310
        43 = 475
    }
    "全部类型" = @{
        4  = 1243
        10 = 3456
        12 = 85
        17 = 592
        18 = 86
        19 = 730
        20 = 207
        25 = 2602
        26 = 5107
        30 = 1534
        32 = 2239
        36 = 113
        38 = 310
        44 = 475
    }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]

    foreach ($row in $rowsForSheet.Keys) {
        $newValue = $rowsForSheet[$row]
        $ws.Cells.Item($row, 6).Value = $newValue
    }
}
